$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.1225990960284
$ws.Range("C2").Value = -7.750390181979499
$ws.Range("D2").Value = 0.4625532265844802
$ws.Range("E2").Value = 1.114035804318275
$ws.Range("F2").Value = -0.6809154673822229
$ws.Range("G2").Value = -1.636980281981076
$ws.Range("H2").Value = 1.198145430806236
$ws.Range("I2").Value = -0.4702713277968628
$ws.Range("J2").Value = 0.550237693327104
$ws.Range("K2").Value = -1.163387135999541

$ws.Range("B3").Value = -8.665941479929575
$ws.Range("C3").Value = -0.09288568806548997
$ws.Range("D3").Value = 0.79511609412051
$ws.Range("E3").Value = -0.8901178518498177
$ws.Range("F3").Value = -1.790197713789883
$ws.Range("G3").Value = 1.072690561200211
$ws.Range("H3").Value = -0.5818437195903187
$ws.Range("I3").Value = 0.4455899393931738
$ws.Range("J3").Value = -1.264578414563814
$ws.Range("K3").Value = 1.565582411450561

$ws.Range("B4").Value = -7.258115570652438
$ws.Range("C4").Value = -4.979538703426429
$ws.Range("D4").Value = -5.532175812781836
$ws.Range("E4").Value = -5.540725542422027
$ws.Range("F4").Value = -1.983394607614172
$ws.Range("G4").Value = -3.101548490755464
$ws.Range("H4").Value = -1.661888171553665
$ws.Range("I4").Value = -3.056284156481218
$ws.Range("J4").Value = 0.01525390354836409
$ws.Range("K4").Value = -2.034922066113249

$ws.Range("B5").Value = -3.299357556857204
$ws.Range("C5").Value = 1.883868949726868
$ws.Range("D5").Value = -3.292669226196934
$ws.Range("E5").Value = 1.997966407575002
$ws.Range("F5").Value = -1.123366588845489
$ws.Range("G5").Value = 0.7501180745166554
$ws.Range("H5").Value = -1.466273777980254
$ws.Range("I5").Value = 1.657826449858163
$ws.Range("J5").Value = -0.7515994264059269
$ws.Range("K5").Value = 0.7022009457793779

$ws.Range("B6").Value = -1.843337576991127
$ws.Range("C6").Value = -1.429798571433311
$ws.Range("D6").Value = 1.237057483404022
$ws.Range("E6").Value = -0.5731355932232373
$ws.Range("F6").Value = 0.5028648335934157
$ws.Range("G6").Value = -1.205689637359233
$ws.Range("H6").Value = 1.614045429628475
$ws.Range("I6").Value = -0.6199119289367258
$ws.Range("J6").Value = 0.73223086762307
$ws.Range("K6").Value = 0.1104179136903787

$ws.Range("B7").Value = -1.021566117811541
$ws.Range("C7").Value = 1.488680527290416
$ws.Range("D7").Value = -0.6550416470463725
$ws.Range("E7").Value = 0.5779568829747013
$ws.Range("F7").Value = -1.099782476845509
$ws.Range("G7").Value = 1.67246207578713
$ws.Range("H7").Value = -0.5534732451066342
$ws.Range("I7").Value = 0.8069128262339065
$ws.Range("J7").Value = 0.1807312055913575
$ws.Range("K7").Value = 0.6808125494453288

$ws.Range("B8").Value = 1.383181315441601
$ws.Range("C8").Value = -0.7811469310117304
$ws.Range("D8").Value = 0.746204951607933
$ws.Range("E8").Value = -1.039118384241627
$ws.Range("F8").Value = 1.683594915161763
$ws.Range("G8").Value = -0.4906104937797037
$ws.Range("H8").Value = 0.8642643540276498
$ws.Range("I8").Value = 0.2252363825957365
$ws.Range("J8").Value = 0.7318864204717475
$ws.Range("K8").Value = 0.06620091605482592

$ws.Range("B9").Value = -1.385587550859469
$ws.Range("C9").Value = 0.6820015742532396
$ws.Range("D9").Value = -0.7475710761844102
$ws.Range("E9").Value = 1.645888039075829
$ws.Range("F9").Value = -0.4991945138531648
$ws.Range("G9").Value = 0.948005431625018
$ws.Range("H9").Value = 0.2600635275237815
$ws.Range("I9").Value = 0.7584508035240445
$ws.Range("J9").Value = 0.1116997568014248
$ws.Range("K9").Value = 0.5327578324921759

$ws.Range("B10").Value = 0.5943607454143283
$ws.Range("C10").Value = -0.7999902782140788
$ws.Range("D10").Value = 1.704534436060835
$ws.Range("E10").Value = -0.4942653685160635
$ws.Range("F10").Value = 0.9369121485761228
$ws.Range("G10").Value = 0.272830226748868
$ws.Range("H10").Value = 0.7671803625714035
$ws.Range("I10").Value = 0.1142017009502766
$ws.Range("J10").Value = 0.539033769963857
$ws.Range("K10").Value = 0.5940197511622507

$ws.Range("B11").Value = -0.7775622985810702
$ws.Range("C11").Value = 1.707137975623284
$ws.Range("D11").Value = -0.5163330541700787
$ws.Range("E11").Value = 0.9309301682717679
$ws.Range("F11").Value = 0.2674869765664569
$ws.Range("G11").Value = 0.7568118292010373
$ws.Range("H11").Value = 0.1057243788266528
$ws.Range("I11").Value = 0.5312614616570462
$ws.Range("J11").Value = 0.5853957776292821
$ws.Range("K11").Value = 0.7117978212943072

$ws.Range("B12").Value = 1.728562547080504
$ws.Range("C12").Value = -0.3961552053708663
$ws.Range("D12").Value = 0.8417358238579847
$ws.Range("E12").Value = 0.2377599727178791
$ws.Range("F12").Value = 0.7654379492309473
$ws.Range("G12").Value = 0.07978420540121761
$ws.Range("H12").Value = 0.5084828986288187
$ws.Range("I12").Value = 0.5717777424806643
$ws.Range("J12").Value = 0.6932770241091315
$ws.Range("K12").Value = -0.2007560336349775

$ws.Range("B13").Value = -0.4384758376912558
$ws.Range("C13").Value = 0.8125313500022515
$ws.Range("D13").Value = 0.2347393729129579
$ws.Range("E13").Value = 0.7465766656194669
$ws.Range("F13").Value = 0.05945307433601332
$ws.Range("G13").Value = 0.4939527431945709
$ws.Range("H13").Value = 0.5551366969883963
$ws.Range("I13").Value = 0.6756931786803813
$ws.Range("J13").Value = -0.2172299171650203
$ws.Range("K13").Value = 0.553751414566908

$ws.Range("B14").Value = 1.152729090620161
$ws.Range("C14").Value = 0.310110434696895
$ws.Range("D14").Value = 0.5554833793064679
$ws.Range("E14").Value = 0.08414698692274653
$ws.Range("F14").Value = 0.5033321467508669
$ws.Range("G14").Value = 0.5003258983770053
$ws.Range("H14").Value = 0.6567990877917116
$ws.Range("I14").Value = -0.2320939179800661
$ws.Range("J14").Value = 0.5250397039375373
$ws.Range("K14").Value = 0.2638965897873631

$ws.Range("B15").Value = 0.7608053066871455
$ws.Range("C15").Value = 0.6025940815188262
$ws.Range("D15").Value = -0.1561362438680617
$ws.Range("E15").Value = 0.5392772278167297
$ws.Range("F15").Value = 0.4982751409374124
$ws.Range("G15").Value = 0.5802308483556813
$ws.Range("H15").Value = -0.2594933123553494
$ws.Range("I15").Value = 0.4987506537398119
$ws.Range("J15").Value = 0.2204992990740305

$ws.Range("B16").Value = 0.9149196684423646
$ws.Range("C16").Value = -0.02284750413253739
$ws.Range("D16").Value = 0.3601988060005381
$ws.Range("E16").Value = 0.5264693797079796
$ws.Range("F16").Value = 0.6162561595480749
$ws.Range("G16").Value = -0.2970348825595631
$ws.Range("H16").Value = 0.4932083146524507
$ws.Range("I16").Value = 0.2246746280127792

$ws.Range("B17").Value = 0.2127429869753038
$ws.Range("C17").Value = 0.4448775468748477
$ws.Range("D17").Value = 0.3808690130742625
$ws.Range("E17").Value = 0.6262099024073174
$ws.Range("F17").Value = -0.2816697128600181
$ws.Range("G17").Value = 0.4543599796950276
$ws.Range("H17").Value = 0.2088288189855932

$ws.Range("B18").Value = 0.7559017333562305
$ws.Range("C18").Value = 0.4979577874854577
$ws.Range("D18").Value = 0.4637457609577506
$ws.Range("E18").Value = -0.2527107623948165
$ws.Range("F18").Value = 0.4901028521499312
$ws.Range("G18").Value = 0.1775011726019661

$ws.Range("B19").Value = 0.7447829648895721
$ws.Range("C19").Value = 0.4814016284956401
$ws.Range("D19").Value = -0.347333001505811
$ws.Range("E19").Value = 0.5223202403984113
$ws.Range("F19").Value = 0.1898892984296834

$ws.Range("B20").Value = 0.7210779879118521
$ws.Range("C20").Value = -0.2623087580365975
$ws.Range("D20").Value = 0.4058356620403972
$ws.Range("E20").Value = 0.2049945700815359

$ws.Range("B21").Value = -0.09744868100251025
$ws.Range("C21").Value = 0.4193729342883134
$ws.Range("D21").Value = 0.1420216510915729

$ws.Range("B22").Value = 0.6732219761537215
$ws.Range("C22").Value = 0.2413397012736094

$ws.Range("B23").Value = 0.2853993925130583


$ws.Range("K15").ClearContents()
$ws.Range("J16").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("H18").ClearContents()
$ws.Range("G19").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("B24").ClearContents()
